# Format Upload Cakupan Imunisasi - rubah format impor (#47)
# Change desa_id column from numeric village codes to string-formatted codes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "53.06.13.2021",
    "53.06.13.2020",
    "53.06.13.2019",
    "53.06.13.2018",
    "53.06.13.2017",
    "53.06.13.2016",
    "53.06.13.2015",
    "53.06.13.2014"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Widen column A to fit the new text values, and drop the old bestFit state
$ws.Columns.Item(1).ColumnWidth = 13.6

# Move the active selection like in the edited workbook
$ws.Range("D9").Select()
